$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D, shifting existing D:K data to F:M
$ws.Range("D:E").Insert(-4161)

# Copy formatting from the old D:E columns (now shifted to F:G) onto the new D:E columns
# (only for the rows that actually contain data cells in D:M)
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new columns (D, E) with the newest-quarter data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 141000
$ws.Range("E8").Value = 130600
$ws.Range("D9").Value = 60800
$ws.Range("E9").Value = 53500
$ws.Range("D10").Value = 80200
$ws.Range("E10").Value = 77100
$ws.Range("D12").Value = 15300
$ws.Range("E12").Value = 15100
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 23000
$ws.Range("E14").Value = 11400
$ws.Range("D15").Value = 9200
$ws.Range("E15").Value = 4500
$ws.Range("D17").Value = 156100
$ws.Range("E17").Value = 133500
$ws.Range("D18").Value = -15100
$ws.Range("E18").Value = -2900
$ws.Range("D20").Value = -1200
$ws.Range("E20").Value = 1000
$ws.Range("D21").Value = -8100
$ws.Range("E21").Value = 7000
$ws.Range("D22").Value = 1500
$ws.Range("E22").Value = 1600
$ws.Range("D23").Value = -17900
$ws.Range("E23").Value = -3600
$ws.Range("D24").Value = -14900
$ws.Range("E24").Value = 1900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -3000
$ws.Range("E26").Value = -5600
$ws.Range("D27").Value = -3000
$ws.Range("E27").Value = -5600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -8600
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 1200
$ws.Range("E32").Value = -1000
$ws.Range("D33").Value = -11600
$ws.Range("E33").Value = -5600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -11600
$ws.Range("E35").Value = -5600
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 56400
$ws.Range("E41").Value = 54400
$ws.Range("D42").Value = "NA"
$ws.Range("E42").Value = "NA"
$ws.Range("D43").Value = 127000
$ws.Range("E43").Value = 121100
$ws.Range("D44").Value = 79700
$ws.Range("E44").Value = 80600
$ws.Range("D45").Value = 22600
$ws.Range("E45").Value = 30800
$ws.Range("D46").Value = 285800
$ws.Range("E46").Value = 287000
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 22900
$ws.Range("E48").Value = 21600
$ws.Range("D49").Value = 287100
$ws.Range("E49").Value = 318600
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 42400
$ws.Range("E52").Value = 26900
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 638100
$ws.Range("E54").Value = 654000
$ws.Range("D57").Value = 28800
$ws.Range("E57").Value = 20600
$ws.Range("D58").Value = 35000
$ws.Range("E58").Value = 20000
$ws.Range("D59").Value = 69600
$ws.Range("E59").Value = 68800
$ws.Range("D60").Value = 133400
$ws.Range("E60").Value = 109400
$ws.Range("D61").Value = 69500
$ws.Range("E61").Value = 94400
$ws.Range("D62").Value = 36800
$ws.Range("E62").Value = 40200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 239700
$ws.Range("E66").Value = 244100
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 102300
$ws.Range("E72").Value = 113900
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 398400
$ws.Range("E76").Value = 410000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -11600
$ws.Range("E81").Value = -5600
$ws.Range("D83").Value = 8200
$ws.Range("E83").Value = 9000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 14100
$ws.Range("E89").Value = 7600
$ws.Range("D91").Value = -2700
$ws.Range("E91").Value = -1700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -2800
$ws.Range("E94").Value = -2100
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -9100
$ws.Range("E100").Value = -4400
$ws.Range("D101").Value = -300
$ws.Range("E101").Value = -1600
$ws.Range("D102").Value = 1900
$ws.Range("E102").Value = -500

# Row 58 and Row 91 received additional data corrections beyond the simple column shift
$ws.Range("D58").Value = 35000
$ws.Range("E58").Value = 20000
$ws.Range("F58").Value = "NA"
$ws.Range("G58").Value = "NA"
$ws.Range("H58").Value = "NA"
$ws.Range("I58").Value = "NA"
$ws.Range("J58").Value = "NA"
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 0
$ws.Range("D91").Value = -2700
$ws.Range("E91").Value = -1700
$ws.Range("F91").Value = -900
$ws.Range("G91").Value = -2500
$ws.Range("H91").Value = -1300
$ws.Range("I91").Value = -1300
$ws.Range("J91").Value = -500
$ws.Range("K91").Value = -1000
$ws.Range("L91").Value = -1000
$ws.Range("M91").Value = 0
